$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 7, carrying over the formatting (styles) from row 6 ---
$ws.Range("A6:CH6").Copy()
$ws.Range("A7:CH7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(7).RowHeight = 15.75

# --- Populate the new row's values ---
$ws.Range("A7").Value = 42823.302823796301
$ws.Range("B7").Value = 'yuzhu.zhang@outlook.de'
$ws.Range("C7").Value = '601.1 National Way of Life General: Positive'
$ws.Range("D7").Value = '601.1 National Way of Life General: Positive'
$ws.Range("E7").Value = '601.1 National Way of Life General: Positive'
$ws.Range("F7").Value = '601.1 National Way of Life General: Positive'
$ws.Range("G7").Value = '601.1 National Way of Life General: Positive'
$ws.Range("H7").Value = '110 European/LA Integration: Negative'
$ws.Range("I7").Value = '305.1 Political Authority: Party Competence'
$ws.Range("J7").Value = '110 European/LA Integration: Negative'
$ws.Range("K7").Value = '000 No meaningful category applies'
$ws.Range("L7").Value = '000 No meaningful category applies'
$ws.Range("M7").Value = '408 Economic Goals'
$ws.Range("N7").Value = '110 European/LA Integration: Negative'
$ws.Range("O7").Value = '109 Internationalism: Negative'
$ws.Range("P7").Value = '302 Centralisation: Positive'
$ws.Range("Q7").Value = '411 Technology and Infrastructure: Positive'
$ws.Range("R7").Value = '403 Market Regulation: Positive'
$ws.Range("S7").Value = '303 Governmental and Administrative Efficiency: Positive'
$ws.Range("T7").Value = '305.1 Political Authority: Party Competence'
$ws.Range("U7").Value = '305.1 Political Authority: Party Competence'
$ws.Range("V7").Value = '504 Welfare State Expansion'
$ws.Range("W7").Value = '504 Welfare State Expansion'
$ws.Range("X7").Value = '504 Welfare State Expansion'
$ws.Range("Y7").Value = '504 Welfare State Expansion'
$ws.Range("Z7").Value = '503 Equality: Positive'
$ws.Range("AA7").Value = '503 Equality: Positive'
$ws.Range("AB7").Value = '504 Welfare State Expansion'
$ws.Range("AC7").Value = '504 Welfare State Expansion'
$ws.Range("AD7").Value = '705 Minority Groups: Positive'
$ws.Range("AE7").Value = '601.2 Immigration: Negative'
$ws.Range("AF7").Value = '402 Incentives: Positive'
$ws.Range("AG7").Value = '506 Education Expansion'
$ws.Range("AH7").Value = '104 Military: Positive'
$ws.Range("AI7").Value = '504 Welfare State Expansion'
$ws.Range("AJ7").Value = '504 Welfare State Expansion'
$ws.Range("AK7").Value = '504 Welfare State Expansion'
$ws.Range("AL7").Value = '504 Welfare State Expansion'
$ws.Range("AM7").Value = '504 Welfare State Expansion'
$ws.Range("AN7").Value = '000 No meaningful category applies'
$ws.Range("AO7").Value = '305.1 Political Authority: Party Competence'
$ws.Range("AP7").Value = '601.1 National Way of Life General: Positive'
$ws.Range("AQ7").Value = '110 European/LA Integration: Negative'
$ws.Range("AR7").Value = '402 Incentives: Positive'
$ws.Range("AS7").Value = '601.2 Immigration: Negative'
$ws.Range("AT7").Value = '504 Welfare State Expansion'
$ws.Range("AU7").Value = '504 Welfare State Expansion'
$ws.Range("AV7").Value = '504 Welfare State Expansion'
$ws.Range("AW7").Value = '504 Welfare State Expansion'
$ws.Range("AX7").Value = '504 Welfare State Expansion'
$ws.Range("AY7").Value = '704 Middle Class and Professional Groups: Positive'
$ws.Range("AZ7").Value = '503 Equality: Positive'
$ws.Range("BA7").Value = '503 Equality: Positive'
$ws.Range("BB7").Value = '503 Equality: Positive'
$ws.Range("BC7").Value = '503 Equality: Positive'
$ws.Range("BD7").Value = '503 Equality: Positive'
$ws.Range("BE7").Value = '701 Labour Groups: Positive'
$ws.Range("BF7").Value = '503 Equality: Positive'
$ws.Range("BG7").Value = '110 European/LA Integration: Negative'
$ws.Range("BH7").Value = '110 European/LA Integration: Negative'
$ws.Range("BI7").Value = '402 Incentives: Positive'
$ws.Range("BJ7").Value = '406 Protectionism: Positive'
$ws.Range("BK7").Value = '606.1 Civic Mindedness General: Positive'
$ws.Range("BL7").Value = '110 European/LA Integration: Negative'
$ws.Range("BM7").Value = '402 Incentives: Positive'
$ws.Range("BN7").Value = '110 European/LA Integration: Negative'
$ws.Range("BO7").Value = '406 Protectionism: Positive'
$ws.Range("BP7").Value = '305.1 Political Authority: Party Competence'
$ws.Range("BQ7").Value = '305.1 Political Authority: Party Competence'
$ws.Range("BR7").Value = '706 Non-Economic Demographic Groups: Positive'
$ws.Range("BS7").Value = '000 No meaningful category applies'
$ws.Range("BT7").Value = '000 No meaningful category applies'
$ws.Range("BU7").Value = '305.1 Political Authority: Party Competence'
$ws.Range("BV7").Value = '000 No meaningful category applies'
$ws.Range("BW7").Value = '000 No meaningful category applies'
$ws.Range("BX7").Value = '303 Governmental and Administrative Efficiency: Positive'
$ws.Range("BY7").Value = '303 Governmental and Administrative Efficiency: Positive'
$ws.Range("BZ7").Value = '305.1 Political Authority: Party Competence'
$ws.Range("CA7").Value = '303 Governmental and Administrative Efficiency: Positive'
$ws.Range("CB7").Value = '110 European/LA Integration: Negative'
$ws.Range("CC7").Value = '109 Internationalism: Negative'
$ws.Range("CD7").Value = '000 No meaningful category applies'
$ws.Range("CE7").Value = '303 Governmental and Administrative Efficiency: Positive'
$ws.Range("CF7").Value = '303 Governmental and Administrative Efficiency: Positive'
$ws.Range("CG7").Value = '303 Governmental and Administrative Efficiency: Positive'
$ws.Range("CH7").Value = '000 No meaningful category applies'

# --- Page margins: top/bottom set to 0.75in (54pt) ---
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
